$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to stay text so that
# numeric-looking values (e.g. "0.9981", "1.210") are not auto-converted
# to numbers by Excel and keep exact formatting/precision.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.888.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.18%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.621.87"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.36%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9981"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.21%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.79"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.03%  "

# Row 6
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.07%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3926"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.35%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3842"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.76%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9978"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.27%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.374"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.52%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "49.30"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.30%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08454"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.48%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.97"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.41%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.061"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.94%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.557"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.29%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001278"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.27%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.610.17"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.88%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.54"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.59%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06909"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.64%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.10"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.14%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.834"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.75%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.43"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.86%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.884.43"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.18%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.435"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.25%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.889"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.27%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.22"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.83%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "156.77"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.14%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "139.34"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.22%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.212"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -11.30%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.916"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.27%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.473"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.86%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.783.95"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.92%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08049"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.26%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9794"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.39%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02887"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.97%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.593"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.53%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2668"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.13%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09185"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.90%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.32"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.31%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.58"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.00%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.423"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.30%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7506"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.39%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.89"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.60%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6898"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.29%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.469"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.89%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.062"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.53%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9995"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.01%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08258"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.62%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.11"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.79%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.210"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -9.29%  "

